$d = $word.ActiveDocument

$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*V1 trains the classification model first on 248*248*") {
        $para = $p
        break
    }
}

$r = $para.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p w:rsidR="0008338E" w:rsidRDefault="006C785D" w:rsidP="00833E95">' +
       '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
       '<w:r><w:t>V1 trains the c</w:t></w:r>' +
       '<w:r><w:t>lassification model first on 224*224</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> resolution data and then increases the resolution to 448*448 for detection.</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
